$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Three literal text replacements in the existing "source file" headers.
# ---------------------------------------------------------------------------

$old1 = "Files\\2015 Case Study\\Primary Sources_Policy_Strategies\\2015 DOD Cyber Strategy CLEAN - § 1 reference coded [ 0.06% Coverage]"
$new1 = "Files\\2015 Case Study\\CS2_Primary Sources_Policy_Strategies\\2015 DoD Cyber Strategy - § 1 reference coded [ 0.06% Coverage]"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "Files\\2015 Case Study\\Primary Sources_Policy_Strategies\\2015 National Security Strategy CLEAN - § 5 references coded [ 0.51% Coverage]"
$new2 = "Files\\2015 Case Study\\CS2_Primary Sources_Policy_Strategies\\2015 National Security Strategy - § 5 references coded [ 0.51% Coverage]"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

$old3 = "Files\\2015 Case Study\\Primary Sources_Policy_Strategies\\2015 WH Report on Cyber Deterrence Policy Final CLEAN - § 1 reference coded [ 0.18% Coverage]"
$new3 = "Files\\2015 Case Study\\CS2_Primary Sources_Policy_Strategies\\2015 White House Report on Cyber Deterrence Policy - § 1 reference coded [ 0.18% Coverage]"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Append 9 new paragraphs at the end of the document (new Case Study block)
#
# Strategy: first append 9 blank paragraphs (each new paragraph inherits the
# plain / non-highlighted style of the current last paragraph, so they all
# start out plain), then go back and fill in text, applying the "header"
# (gray highlight + hanging indent) formatting only to the 4 header rows.
# This avoids highlight/indent formatting leaking into the plain paragraphs
# that would otherwise inherit it from a preceding header paragraph.
# ---------------------------------------------------------------------------

$lb = [char]0x0B

$newTexts = @(
    "Files\\2018 Case Study\\CS3_Primary Sources_Policy_Strategies\\2017 National Security Strategy - § 4 references coded [ 0.24% Coverage]",
    "Reference 1 - 0.06% Coverage",
    "As a growing supplier of energy resources, technologies, and services around the world, the United States will help our allies and partners become more resilient against those that use energy to coerce.",
    "Reference 2 - 0.05% Coverage",
    "We will not allow adversaries to use threats of nuclear escalation or other irresponsible nuclear behaviors to coerce the United States, our allies, and our partners.",
    "Reference 3 - 0.06% Coverage",
    "Economic tools—including sanctions, anti-money-laundering and anti-corruption measures, and enforcement actions—can be important parts of broader strategies to deter, coerce, and constrain adversaries.",
    "Reference 4 - 0.07% Coverage",
    ("We will maintain our strong ties with Taiwan in accordance with our “One China” policy, including our commitments under the " + $lb + "47 " + $lb + "Taiwan Relations Act to provide for Taiwan’s legitimate defense needs and deter coercion.")
)

# true => this row gets the gray-highlighted "header" paragraph formatting
$isHeader = @($true, $true, $false, $true, $false, $true, $false, $true, $false)

$startCount = $d.Paragraphs.Count

for ($i = 0; $i -lt $newTexts.Count; $i++) {
    $tail = $d.Paragraphs.Last.Range
    $tail.InsertParagraphAfter() | Out-Null
}

for ($i = 0; $i -lt $newTexts.Count; $i++) {
    $p = $d.Paragraphs.Item($startCount + $i + 1)
    $p.Range.Text = $newTexts[$i]
    if ($isHeader[$i]) {
        $p.SpaceBefore = 5.65
        $p.SpaceAfter = 5.65
        $p.LeftIndent = 5.65
        $p.RightIndent = 5.65
        $p.FirstLineIndent = -0.0001
        $p.Range.Font.HighlightColorIndex = 16
    }
}

Write-Output "Done"
